# Generate Report for Handoff
# Adds two new rows (for 0092fc96-548b-4cf1-83ca-c94a2c4ed957 and
# 53bff250-c3e0-4453-992f-81b721c9b706) to the Overview / zh-cn / de-de
# sheets of the localization-status workbook.

$wb = $excel.ActiveWorkbook

$hyperColor = 15570276   # OLE BGR for RGB(0x64,0x95,0xED) == style "FF6495ED"

function Add-LinkCell {
    param($ws, $cellRef, $address, $text)

    $ws.Hyperlinks.Add($ws.Range($cellRef), $address, "", "", $text) | Out-Null
    $ws.Range($cellRef).Font.Name = "Calibri"
    $ws.Range($cellRef).Font.Underline = $true
    $ws.Range($cellRef).Font.Color = $hyperColor
}

# ---------------------------------------------------------------------
# Source-repo commit / hash constants (reusing the same commit shas the
# existing rows already point at).
# ---------------------------------------------------------------------
$mdCommit  = "92bce5e15a93f342cd1d079aa3289bf6d5525d48"
$zhCommit  = "956fef6632102003fa8cff20f08de19b629e0b4a"
$deCommit  = "e82c6989a635019ba0f8bd0209944665f9537821"

$uuid1 = "0092fc96-548b-4cf1-83ca-c94a2c4ed957"
$uuid1Md = "0092fc96-548b-4cf1-83ca-c94a2c4ed957.md"
$uuid1ZhXlf = "0092fc96-548b-4cf1-83ca-c94a2c4ed957.7daadf778786066de6401f8f09645479cd44e3d3.zh-cn.xlf"
$uuid1DeXlf = "0092fc96-548b-4cf1-83ca-c94a2c4ed957.7daadf778786066de6401f8f09645479cd44e3d3.de-de.xlf"

$uuid2 = "53bff250-c3e0-4453-992f-81b721c9b706"
$uuid2Md = "53bff250-c3e0-4453-992f-81b721c9b706.md"
$uuid2ZhXlf = "53bff250-c3e0-4453-992f-81b721c9b706.8ae5f62414bf45c28f3d36ef164dbb033427489f.zh-cn.xlf"
$uuid2DeXlf = "53bff250-c3e0-4453-992f-81b721c9b706.8ae5f62414bf45c28f3d36ef164dbb033427489f.de-de.xlf"

$status = "Ready for handoff"

$mdUrl1 = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$uuid1.md"
$mdUrl2 = "https://github.com/OpenLocalizationTest/oltest/blob/$mdCommit/e2e/$uuid2.md"

$zhUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$uuid1ZhXlf"
$zhUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/$uuid2ZhXlf"

$deUrl1 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$uuid1DeXlf"
$deUrl2 = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/$uuid2DeXlf"

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

Add-LinkCell $wsOverview "A4" $mdUrl1 $uuid1Md
$wsOverview.Range("B4").Value2 = $status
$wsOverview.Range("C4").Value2 = $status
$wsOverview.Range("D4").Value2 = "2016-02-14 06:02:18"

Add-LinkCell $wsOverview "A5" $mdUrl2 $uuid2Md
$wsOverview.Range("B5").Value2 = $status
$wsOverview.Range("C5").Value2 = $status
$wsOverview.Range("D5").Value2 = "2016-02-14 06:02:18"

# ---------------------------------------------------------------------
# Sheet "zh-cn": Source File Name | File Extension | Status |
#   Latest Handoff File | Latest Handoff Datetime | Latest Target File |
#   Latest Handback File | Latest Handback DateTime | Handoff Reason |
#   Dependency From | Error Detail
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

Add-LinkCell $wsZh "A4" $mdUrl1 $uuid1Md
Add-LinkCell $wsZh "B4" $mdUrl1 ".md"
$wsZh.Range("C4").Value2 = $status
Add-LinkCell $wsZh "D4" $zhUrl1 $uuid1ZhXlf
$wsZh.Range("E4").Value2 = "2016-03-14 06:02:10"
$wsZh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H4").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("I4").Value2 = "Include"

Add-LinkCell $wsZh "A5" $mdUrl2 $uuid2Md
Add-LinkCell $wsZh "B5" $mdUrl2 ".md"
$wsZh.Range("C5").Value2 = $status
Add-LinkCell $wsZh "D5" $zhUrl2 $uuid2ZhXlf
$wsZh.Range("E5").Value2 = "2016-03-14 06:02:10"
$wsZh.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsZh.Range("H5").Value2 = "0001-01-01 00:00:00"
$wsZh.Range("I5").Value2 = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de": same columns as zh-cn
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

Add-LinkCell $wsDe "A4" $mdUrl1 $uuid1Md
Add-LinkCell $wsDe "B4" $mdUrl1 ".md"
$wsDe.Range("C4").Value2 = $status
Add-LinkCell $wsDe "D4" $deUrl1 $uuid1DeXlf
$wsDe.Range("E4").Value2 = "2016-03-14 06:02:18"
$wsDe.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H4").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("I4").Value2 = "Include"

Add-LinkCell $wsDe "A5" $mdUrl2 $uuid2Md
Add-LinkCell $wsDe "B5" $mdUrl2 ".md"
$wsDe.Range("C5").Value2 = $status
Add-LinkCell $wsDe "D5" $deUrl2 $uuid2DeXlf
$wsDe.Range("E5").Value2 = "2016-03-14 06:02:18"
$wsDe.Range("E5").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$wsDe.Range("H5").Value2 = "0001-01-01 00:00:00"
$wsDe.Range("I5").Value2 = "Include"

Write-Output "Report rows for handoff generated."
